$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update the "Förändrad" (Changed) date column C for rows 2-17 from
# serial date 45177 (2023-09-08) to 45178 (2023-09-09), keeping existing
# cell formatting/style intact.
foreach ($row in 2..17) {
    $ws.Cells.Item($row, 3).Value = 45178
}
